$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45204 -> 2023-10-05) that
# needs to be bumped by one day (45205 -> 2023-10-06) for every data row
# (rows 2 through 498).
$lastRow = 498
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45205
